$d = $word.ActiveDocument

# Locate the target paragraph: the last paragraph of the document, which
# currently reads "복수 테이블 검색 " - "여러 테이블을 대상" [bookmark] "으로 검색을 수행 "
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*복수 테이블 검색*") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "could not locate '복수 테이블 검색' paragraph"
}

$full = $target.Range

# Replace the whole paragraph's contents (this also removes the pPr rFonts
# override and the _GoBack bookmark, merging the text into 3 clean runs).
$replacementXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">복수 테이블 검색 </w:t></w:r><w:r><w:t xml:space="preserve">– </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">여러 테이블을 대상으로 검색을 수행 </w:t></w:r></w:p>'
$full.InsertXML($replacementXml)

# Re-fetch the (now rewritten) paragraph and insert all of the new content
# that follows it, in one shot, right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*복수 테이블 검색*") {
        $target = $p
    }
}
$insertionPoint = $target.Range.Duplicate
$insertionPoint.Collapse(0)
$insertionPoint.MoveEnd(1, -1)

$newParasXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="22" /><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF" /></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22" /><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t>106 DML – SELECT -2</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>그룹함수</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>– GROUP BY</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>절에 지정된 그룹별로 속성의 값을 집계할 때 사용</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>WINDOW</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">함수 </w:t></w:r><w:r><w:t xml:space="preserve">– </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">GROUP </w:t></w:r><w:r><w:t>BY</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>절을 이용하지 않고 함수의 인수로 지정한 속성의 값을 집계</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">그룹 지정 검색 </w:t></w:r><w:r><w:t>– GROUP BY</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">절에 지정한 속성을 기준으로 자료를 그룹화하여 검색 </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">집합 연산자를 이용한 통합 질의 </w:t></w:r><w:r><w:t xml:space="preserve">– </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">집합 연산자를 사용하여 </w:t></w:r><w:r><w:t>2</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>개 이상의 테이블의 데이터를 하나로 통합</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">UNION / UNION ALL / INTERSECT / EXCEPT </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia" /><w:sz w:val="22" /><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF" /></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22" /><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t xml:space="preserve">107 </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /><w:sz w:val="22" /><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t xml:space="preserve">DML </w:t></w:r><w:r><w:rPr><w:sz w:val="22" /><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t>–</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /><w:sz w:val="22" /><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t xml:space="preserve"> JOIN</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">JOIN </w:t></w:r><w:r><w:t>–</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve"> 2개의 릴레이션에서 연관된 </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>튜플들을</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve"> 결합하여 하나의 새로운 </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>릴레이션을</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve"> 반환</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>INNER JOIN</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">– EQUI </w:t></w:r><w:proofErr w:type="gramStart" /><w:r><w:t>JOIN :</w:t></w:r><w:proofErr w:type="gramEnd" /><w:r><w:t xml:space="preserve"> JOIN </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>대상 테이블에서 공통 속성을 기준으로 비교에 의해 같은 값을 가지는 행을 연결하여 결과를 생성하는 방법</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak /><w:t xml:space="preserve">- </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">NON-EQUI </w:t></w:r><w:proofErr w:type="gramStart" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>JOIN :</w:t></w:r><w:proofErr w:type="gramEnd" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve"> JOIN</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">조건에 </w:t></w:r><w:r><w:t>=</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve"> 조건이 아닌 나머지 비교 연산자를 사용하는 방법</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">OUTER JOIN </w:t></w:r><w:r><w:t>–</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve"> 릴레이션에서 </w:t></w:r><w:r><w:t xml:space="preserve">JOIN </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">조건에 만족하지 않는 </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>튜플도</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve"> 결과로 출력하기 위한 </w:t></w:r><w:r><w:t xml:space="preserve">JOIN </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>방법</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>LEFT OUTER JOIN / RIGHT OUTER JOIN / FULL OUTER JOIN</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack" /><w:bookmarkEnd w:id="0" /></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr></w:pPr></w:p>'
$insertionPoint.InsertXML($newParasXml)

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
